$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold numeric-looking / percentage
# TEXT (e.g. "310.84", "0.60%", "--", "--%") in the source file, so each
# touched D/E cell is forced to a Text format before the write to stop
# Excel re-interpreting the string as a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.02%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.090"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.80%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07571"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.19%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.299"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.657"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.17%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9291"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.15%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1215"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.26%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1797"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.76%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04154"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.22%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.26%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001292"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.56%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005768"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.53%"
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.004085"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.30%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.350"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3354"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.55%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.649"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.70%"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1351"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.39%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2809"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.93%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04023"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.36%"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001267"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.61%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001272"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.04%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02428"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.96%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05147"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.02%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007742"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.21%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.22%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007655"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "12.10%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "14.73%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007992"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.43%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3107"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.32%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006584"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.04%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2703"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "31.54%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "2.48%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.05%"
